# Atualizando o arquivo XLSX
# Apply updated odds values to Sheet1 as described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 2.63
$ws.Range("I2").Value = 2.75
$ws.Range("N2").Value = 8.5
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 10
$ws.Range("AK2").Value = 29
$ws.Range("AW2").Value = 4.75

# Row 11
$ws.Range("G11").Value = 4.75
$ws.Range("H11").Value = 3.55
$ws.Range("J11").Value = 4.9
$ws.Range("O11").Value = 1.33
$ws.Range("P11").Value = 3.05
$ws.Range("W11").Value = 12
$ws.Range("AE11").Value = 16.5
$ws.Range("AG11").Value = 700
$ws.Range("AH11").Value = 6.3
$ws.Range("AI11").Value = 7.5
$ws.Range("AM11").Value = 29
$ws.Range("AQ11").Value = 150
$ws.Range("AS11").Value = 400
$ws.Range("AU11").Value = 7.5
$ws.Range("AW11").Value = 3.55

# Row 16
$ws.Range("S16").Value = 1.44
$ws.Range("T16").Value = 2.63

# Row 17
$ws.Range("S17").Value = 1.33
